# Automatische test-sync: 2025-07-31 21:52:50
# Append new log row (#16) to the "Logs" sheet and refresh the dependent
# "Dashboard" summary sheet + conditional-formatting ranges.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$newRow = 16

$logs.Cells.Item($newRow, 1).Value2 = "Heb je de CE-certificaten van dit product?"
$logs.Cells.Item($newRow, 2).Value2 = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value2 = "Testmail #14: Heb je de CE-certificaten van dit product?"
$logs.Cells.Item($newRow, 4).Value2 = "Productinformatie"
$logs.Cells.Item($newRow, 5).Value2 = "Geachte klant,`r`nDank u voor uw e-mail. Voor het verkrijgen van de CE-certificaten van het product waar u naar vraagt, verzoeken wij u ons het specifieke productnummer of de productnaam te verstrekken. Met deze informatie kunnen wij u de relevante certificaten verstrekken.`r`nAls u verdere vragen heeft of meer ondersteuning nodig heeft, aarzel dan niet om contact met ons op te nemen.`r`nMet vriendelijke groet,`r`n[Naam] E-mailassistent - [Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value2 = "2025-07-31 21:52:33"
$logs.Cells.Item($newRow, 7).Value2 = "Ja"
$logs.Cells.Item($newRow, 8).Value2 = "Nee"
$logs.Cells.Item($newRow, 9).Value2 = "Ja"
$logs.Cells.Item($newRow, 10).Value2 = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) so they keep covering
# the whole table down through the newly added row.
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "15")
    $newRange = $logs.Range($col + "2:" + $col + "16")
    $rules = $oldRange.FormatConditions
    if ($rules.Count -gt 0) {
        $rules.Item(1).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary: "Productinformatie" count goes from 4 to 5.
$dash.Range("B3").Value2 = 5
